# Sprint 6, 7, 8 Stories content update
# Appends 38 new fitchFieldId values to Sheet1 (rows 3590-3627)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @(
    "FC_PERIOD_DT_FIR",
    "FC_STATEMENT_ID_FIR",
    "FC_TOTAL_ASSETS_FIR",
    "FC_COUNTRY_RISK_IND_FIR",
    "FC_REGION_FIR",
    "FC_PERIOD_DT_RANK_FIR",
    "FC_PROFIT_FIR",
    "FC_LOAN_QUAL_FIR",
    "FC_MODEL_SCORE_FIR",
    "FC_FIR",
    "FC_BAND_RANK_FIR",
    "FC_PROFIT_CUTOFF_LOW_FIR",
    "FC_PROFIT_CUTOFF_HIGH_FIR",
    "FC_LOAN_QUAL_CUTOFF_LOW_FIR",
    "FC_LOAN_QUAL_CUTOFF_HIGH_FIR",
    "FC_TOTAL_ASSETS_NORM_MEAN_FIR",
    "FC_TOTAL_ASSETS_NORM_SD_FIR",
    "FC_PROFIT_NORM_MEAN_FIR",
    "FC_PROFIT_NORM_SD_FIR",
    "FC_LOAN_QUAL_NORM_MEAN_FIR",
    "FC_LOAN_QUAL_NORM_SD_FIR",
    "FC_PROFIT_COEFF_FIR",
    "FC_LOAN_QUAL_COEFF_FIR",
    "FC_TOTAL_ASSETS_COEFF_FIR",
    "FC_CRI_MODEL_COEFF_FIR",
    "FC_INTERCEPTS_NO_FIR",
    "FC_INTERCEPTS_FIR",
    "FC_MIN_MODEL_SCORE_FIR",
    "FC_MAX_MODEL_SCORE_FIR",
    "FC_NOTCH_DIFF_ALL_FIR",
    "FC_ENTITIES_NOTCH_DIFF_ALL_FIR",
    "FC_NOTCH_DIFF_CF_ALL_FIR",
    "FC_NOTCH_DIFF_EM_FIR",
    "FC_ENTITIES_NOTCH_DIFF_EM_FIR",
    "FC_NOTCH_DIFF_CF_EM_FIR",
    "FC_NOTCH_DIFF_DM_FIR",
    "FC_ENTITIES_NOTCH_DIFF_DM_FIR",
    "FC_NOTCH_DIFF_CF_DM_FIR"
)

$startRow = 3590
$row = $startRow
foreach ($val in $values) {
    $ws.Cells.Item($row, 1).Value = $val
    $row = $row + 1
}

# Match the workbook's final view state: scrolled down, single-cell selection on D3599
[void]$ws.Range("D3599").Select()

Write-Host "Added" $values.Count "fitchFieldId rows starting at row" $startRow
